$wb = $excel.ActiveWorkbook

# --- 1. Rename the icon file referenced in the "icons" sheet ---------------
# BCN_Logo3.png -> BCN_Logo.png  (shared string used by cells B2 and B3)
$wsIcons = $wb.Worksheets.Item("icons")
$wsIcons.Range("B2").Value = "BCN_Logo.png"
$wsIcons.Range("B3").Value = "BCN_Logo.png"

# --- 2. Restore / update the per-sheet cursor position ----------------------
# The other sheets keep the same active cell they already had; only the
# "icons" sheet (the active tab) gets a new selected range.
$wsFieldnames = $wb.Worksheets.Item("fieldnames")
$wsFieldnames.Select() | Out-Null
$wsFieldnames.Range("E47").Select() | Out-Null

$wsUrl = $wb.Worksheets.Item("URL")
$wsUrl.Select() | Out-Null
$wsUrl.Range("B8").Select() | Out-Null

$wsColor = $wb.Worksheets.Item("color")
$wsColor.Select() | Out-Null
$wsColor.Range("A3").Select() | Out-Null

$wsComments = $wb.Worksheets.Item("comments")
$wsComments.Select() | Out-Null
$wsComments.Range("B2").Select() | Out-Null

# "icons" is the active sheet; select B10:B11 there (B10 becomes active cell)
$wsIcons.Select() | Out-Null
$wsIcons.Range("B10:B11").Select() | Out-Null
